$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text values that look numeric (e.g. "246.40").
# Force text formatting first so Excel does not coerce them into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.488.53'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.728.55'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '246.40'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4828'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = '0.2663'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').Value = '0.06221'
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '1.729.91'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = '0.07072'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '15.60'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').Value = '4.593'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').Value = '0.6098'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('D15').Value = '77.27'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '26.489.25'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '0.000007194'
$ws.Range('E19').Value = '  +4.39%  '
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').Value = '1.952.51'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').Value = '4.500'
$ws.Range('E22').Value = '  -2.97%  '
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('D24').Value = '5.242'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('D25').Value = '137.67'
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('D26').Value = '15.41'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '1.775'
$ws.Range('E27').Value = '  -2.03%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').Value = '108.11'
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '1.397'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('D30').Value = '3.975'
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('D31').Value = '0.07985'
$ws.Range('E31').Value = '  +1.37%  '
$ws.Range('D32').Value = '3.689'
$ws.Range('E32').Value = '  -1.60%  '
$ws.Range('D33').Value = '0.04566'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').Value = '0.6339'
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('D37').Value = '0.8896'
$ws.Range('E37').Value = '  -5.64%  '
$ws.Range('D38').Value = '2.020'
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('D39').Value = '2.389'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').Value = '0.01501'
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('E42').Value = '  -10.33%  '
$ws.Range('D43').Value = '5.490'
$ws.Range('E43').Value = '  -4.96%  '
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('D45').Value = '7.004'
$ws.Range('E45').Value = '  +4.00%  '
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('D47').Value = '0.05383'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('D48').Value = '7.891'
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').Value = '30.54'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('D50').Value = '1.251'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('D51').Value = '51.46'
$ws.Range('E51').Value = '  -0.51%  '
